$d = $word.ActiveDocument

# --- First paragraph: add a paragraph border, widen the left indent ---
$p = $d.Paragraphs(1)
$p.Format.Borders.DistanceFromTop = 5
$p.Format.Borders.DistanceFromLeft = 5
$p.Format.Borders.DistanceFromBottom = 5
$p.Format.Borders.DistanceFromRight = 5
$p.Format.LeftIndent = 11.25

# --- Replace the topic-id placeholder text and drop the trailing run ---
# (the trailing " " run is swallowed by extending the search text over it,
# which merges what remains back into a single run with no preserved space)
$d.Content.Find.Execute("**ID__AFFARS_mp_5325_7003_3_topic_4__ID** ", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_MP5325_7003_3_3__ID**", 2)
